# Update "想去人数" (F column) counts on the 展览 and 全部类型 sheets.
# Both sheets carry duplicate data and need the identical set of updates.

$wb = $excel.ActiveWorkbook

$updates = @{
    4  = 77
    5  = 3083
    7  = 2506
    11 = 1269
    15 = 1131
    16 = 316
    21 = 68
    23 = 230
    25 = 256
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
